$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: CheckID 13 - Log File Growing ---
$ws.Range("A17").Value = 13
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "SQL Server Internal Maintenance"
$ws.Range("D17").Value = "Log File Growing"
$ws.Range("E17").Value = "http://BrentOzar.com/go/logsize"

# --- Row 18: CheckID 14 - Log File Shrinking ---
$ws.Range("A18").Value = 14
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = "SQL Server Internal Maintenance"
$ws.Range("D18").Value = "Log File Shrinking"
$ws.Range("E18").Value = "http://BrentOzar.com/go/logsize"

# --- Row 19: CheckID 15 - Compilations/Sec High ---
$ws.Range("A19").Value = 15
$ws.Range("B19").Value = 50
$ws.Range("C19").Value = "Query Problems"
$ws.Range("D19").Value = "Compilations/Sec High"

# --- Row 20: CheckID 16 - Re-Compilations/Sec High ---
$ws.Range("A20").Value = 16
$ws.Range("B20").Value = 50
$ws.Range("C20").Value = "Query Problems"
$ws.Range("E20").Value = "http://BrentOzar.com/go/recompile"
$ws.Range("D20").Value = "Re-Compilations/Sec High"

# URL for row 19 entered last (matches original authoring order)
$ws.Range("E19").Value = "http://BrentOzar.com/go/compile"

# --- Hyperlinks for the new URL cells ---
$ws.Hyperlinks.Add($ws.Range("E17"), "http://BrentOzar.com/go/logsize")
$ws.Range("E17").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E18"), "http://BrentOzar.com/go/logsize")
$ws.Range("E18").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E19"), "http://BrentOzar.com/go/compile")
$ws.Range("E19").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("E20"), "http://BrentOzar.com/go/recompile")
$ws.Range("E20").Style = "Hyperlink"

# --- Update selection to reflect the new last row ---
$ws.Range("A21").Select()
